$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Betarraga" at Macroferia Regional
# de Talca. It belongs chronologically before the existing row 165, so the
# row is inserted there and the remaining records (old rows 165-262) shift
# down by one (to 166-263) automatically.
$ws.Rows(165).Insert()

$ws.Range("A165").Value = 5
$ws.Range("B165").Value = "Macroferia Regional de Talca"
$ws.Range("C165").Value = "Maule"
$ws.Range("D165").Value = 44606
$ws.Range("E165").Value = 7
$ws.Range("F165").Value = 100114014
$ws.Range("G165").Value = "Betarraga"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Segunda"
$ws.Range("J165").Value = 3000
$ws.Range("K165").Value = 700
$ws.Range("L165").Value = 700
$ws.Range("M165").Value = 700
$ws.Range("N165").Value = "$/paquete 5 unidades"
$ws.Range("O165").Value = "Región del Maule"
$ws.Range("P165").Value = 140
$ws.Range("Q165").Value = 5
$ws.Range("R165").Value = "Hortaliza"
